$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = 0.86494175579651977
$ws.Range("D2").Value = 0.71113969518011388
$ws.Range("S2").Value = 0.99253098778207971
$ws.Range("BP2").Value = 0.95082890971610134
$ws.Range("A3").Value = 0.97804617677218997
$ws.Range("E3").Value = 0.9511598984549996
$ws.Range("BC3").Value = 0.89831004232798572
$ws.Range("E4").Value = 0.98085641493039422
$ws.Range("G5").Value = 0.93527419408008761
$ws.Range("U5").Value = 0.9090899307307988
$ws.Range("D6").Value = 0.77352971311139374
$ws.Range("K6").Value = 0.862560269899039
$ws.Range("AP6").Value = 0.8537389596014261
$ws.Range("H7").Value = 0.79376967968468615
$ws.Range("I7").Value = 0.57927850371317935
$ws.Range("I8").Value = 0.95401152011504442
$ws.Range("P8").Value = 0.51514355643695475
$ws.Range("AF8").Value = 0.72494190745062259
$ws.Range("K9").Value = 0.84515130778483694
$ws.Range("R9").Value = 0.93912943551649275
$ws.Range("K10").Value = 0.9809934137143097
$ws.Range("L11").Value = 0.879926316351756
$ws.Range("J12").Value = 0.82535162005938167
$ws.Range("W12").Value = 0.77487531980408786
$ws.Range("AL12").Value = 0.7980285076674869
$ws.Range("L13").Value = 0.95426206835939187
$ws.Range("N13").Value = 0.98561633416525019
$ws.Range("O13").Value = 0.94217000616184698
$ws.Range("P14").Value = 0.97077368072865844
$ws.Range("K15").Value = 0.92473298387252889
$ws.Range("N15").Value = 0.70879074095065242
$ws.Range("P15").Value = 0.82688411510357929
$ws.Range("O17").Value = 0.97031073079799501
$ws.Range("P17").Value = 0.97565376087478239
$ws.Range("R17").Value = 0.90799387477941096
$ws.Range("T18").Value = 0.90344653800042296
$ws.Range("Q19").Value = 0.76321265995256704
$ws.Range("R19").Value = 0.99517141390823172
$ws.Range("AR19").Value = 0.9964696141621181
$ws.Range("B20").Value = 0.65471252618155651
$ws.Range("V20").Value = 0.75433605053535069
$ws.Range("S21").Value = 0.6352798094763441
$ws.Range("V21").Value = 0.79335778196566764
$ws.Range("W21").Value = 0.73270226618308576
$ws.Range("W22").Value = 0.92954011992075003
$ws.Range("O23").Value = 0.96834854770808465
$ws.Range("Y23").Value = 0.5491500946186243
$ws.Range("Y24").Value = 0.91062430374313075
$ws.Range("AN24").Value = 0.91682618515346925
$ws.Range("Y26").Value = 0.86797774968833041
$ws.Range("AA26").Value = 0.81079982260070738
$ws.Range("AB26").Value = 0.95552790569773804
$ws.Range("Y27").Value = 0.88567925521079416
$ws.Range("AB27").Value = 0.81208975706717579
$ws.Range("D28").Value = 0.7320776567980456
$ws.Range("AA29").Value = 0.8151553659167281
$ws.Range("AB29").Value = 0.75331021378953578
$ws.Range("AL29").Value = 0.68759339061360669
$ws.Range("AZ29").Value = 0.93619684057897268
$ws.Range("BP29").Value = 0.87920749520032904
$ws.Range("V30").Value = 0.87367197738277103
$ws.Range("BG30").Value = 0.79319935397902364
$ws.Range("AD31").Value = 0.75003823330987085
$ws.Range("AN31").Value = 0.9279506040134502
$ws.Range("AE33").Value = 0.97808271968892546
$ws.Range("AF33").Value = 0.75713320143381657
$ws.Range("AI33").Value = 0.82241966452964588
$ws.Range("AF34").Value = 0.81623816461955578
$ws.Range("AG34").Value = 0.96423096691687593
$ws.Range("AI34").Value = 0.97906895901752611
$ws.Range("BI34").Value = 0.83741169057015719
$ws.Range("AK35").Value = 0.90888306912029093
$ws.Range("X36").Value = 0.94815173833801891
$ws.Range("AH36").Value = 0.71778020973508638
$ws.Range("AJ37").Value = 0.76124173568673237
$ws.Range("AL37").Value = 0.99580223777855537
$ws.Range("AM37").Value = 0.93982963220494242
$ws.Range("AJ38").Value = 0.80946924068865034
$ws.Range("AN39").Value = 0.97296947017021596
$ws.Range("AO39").Value = 0.93788650669777485
$ws.Range("AL40").Value = 0.96964004541771287
$ws.Range("I41").Value = 0.99449580294329065
$ws.Range("AO42").Value = 0.95843515573915661
$ws.Range("AQ42").Value = 0.86264868485706914
$ws.Range("AR42").Value = 0.71326763454330999
$ws.Range("X43").Value = 0.81422191838271196
$ws.Range("AA43").Value = 0.88377293103978705
$ws.Range("AY43").Value = 0.96694751238115217
$ws.Range("AI44").Value = 0.99379535953257681
$ws.Range("AT45").Value = 0.78638783442001914
$ws.Range("AU45").Value = 0.87745757848988393
$ws.Range("AR46").Value = 0.9589937963867402
$ws.Range("AY46").Value = 0.92650888757815264
$ws.Range("AT47").Value = 0.93068783117341491
$ws.Range("BI47").Value = 0.97645134990317728
$ws.Range("AT48").Value = 0.80733903783339556
$ws.Range("AU48").Value = 0.98204331601003747
$ws.Range("AU49").Value = 0.84795226867416651
$ws.Range("AV49").Value = 0.86507275018628982
$ws.Range("AV50").Value = 0.52184456375074206
$ws.Range("AW50").Value = 0.68882380982338409
$ws.Range("BE50").Value = 0.83894362074360485
$ws.Range("AN51").Value = 0.61180049166992212
$ws.Range("AX52").Value = 0.6298027276667979
$ws.Range("BA52").Value = 0.73906177738324086
$ws.Range("BC53").Value = 0.78744830367136587
$ws.Range("BF53").Value = 0.99345372820961642
$ws.Range("BA54").Value = 0.91291409252612921
$ws.Range("BC54").Value = 0.67369016524518266
$ws.Range("X55").Value = 0.87925690844552418
$ws.Range("BB56").Value = 0.83230740669370529
$ws.Range("BC56").Value = 0.95742433663168758
$ws.Range("BD57").Value = 0.75168035527736987
$ws.Range("BF57").Value = 0.73701473097783887
$ws.Range("BD58").Value = 0.96380409387388033
$ws.Range("BE59").Value = 0.63861972624823027
$ws.Range("BF59").Value = 0.81942376241885317
$ws.Range("BI59").Value = 0.86034665283739575
$ws.Range("BK59").Value = 0.78953188887732961
$ws.Range("BG60").Value = 0.96630204573285083
$ws.Range("BI60").Value = 0.75584747209910619
$ws.Range("BH62").Value = 0.71311015620947049
$ws.Range("BI62").Value = 0.78243795392675064
$ws.Range("BK62").Value = 0.67138134975780894
$ws.Range("BN62").Value = 0.70027770675304457
$ws.Range("AL63").Value = 0.94661957007123454
$ws.Range("A64").Value = 0.75034048365519812
$ws.Range("BM64").Value = 0.60005744243942361
$ws.Range("BN64").Value = 0.74587599146469308
$ws.Range("AS65").Value = 0.69876880792197049
$ws.Range("BN65").Value = 0.78490445288680299
$ws.Range("BO66").Value = 0.56969699250905814
$ws.Range("BP66").Value = 0.98210126667204767
$ws.Range("A67").Value = 0.85492307212591745
$ws.Range("A68").Value = 0.91510867223933567
$ws.Range("BO68").Value = 0.92446215865320347
